$wb = $excel.ActiveWorkbook
$wsSignups = $wb.Worksheets.Item("Signups")
$wsLog = $wb.Worksheets.Item("Log")

# --- Signups sheet: restore (un-hide) the soft-delete column and fix its header text ---
$wsSignups.Columns.Item(1).Hidden = $false
$wsSignups.Range("A1").Value2 = "Deleted"
$wsSignups.Range("A1").Select() | Out-Null

# --- Log sheet: insert a matching "Deleted" column at the front ---
$wsLog.Columns.Item(1).Insert()
$wsLog.Columns.Item(1).ColumnWidth = 12.44140625

# Copy the header style/format from Signups!A1 onto the new Log!A1 header cell
$wsSignups.Range("A1").Copy()
$wsLog.Range("A1").PasteSpecial(-4122)  # xlPasteFormats
$wsLog.Range("A1").Value2 = "Deleted"

# Fill the data rows with the (text) "0" soft-delete flag, keeping default formatting
$dataRange = $wsLog.Range("A2:A5")
$dataRange.NumberFormat = "@"
$wsLog.Cells.Item(2, 1).Value2 = "0"
$wsLog.Cells.Item(3, 1).Value2 = "0"
$wsLog.Cells.Item(4, 1).Value2 = "0"
$wsLog.Cells.Item(5, 1).Value2 = "0"
$dataRange.Style = "Normal"

# Reset the autofilter to the shifted header range
$wsLog.AutoFilterMode = $false
$wsLog.Range("B1:F1").AutoFilter() | Out-Null

# Update the _FilterDatabase defined name to match the new autofilter range
foreach ($n in $wb.Names) {
    if ($n.Name -eq "Log!_FilterDatabase") {
        $n.RefersTo = "=Log!`$B`$1:`$F`$1"
    }
}
